$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New remark text (new shared string) for row 31, column D
$newText = "Makefile: Dependency files incorporated, Manual continued"

# Append a new row of data (row 31): Date, Effort[h], Additional Effort[h], Remark
$ws.Cells.Item(31, 1).Value = 41207
$ws.Cells.Item(31, 1).NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Cells.Item(31, 2).Value = 2.5
$ws.Cells.Item(31, 3).Value = 1.25
$ws.Cells.Item(31, 4).Value = $newText

# Update selection / active cell to reflect the new state (A32 selected after the new row)
$ws.Range("A32").Select()
